# Daily attendance processing - 2025-11-08 07:41:50
#
# The "Recorded By" column (G) lists the users who recorded/edited a
# session, as a comma-separated string. This pass rotates that list for
# every row whose first recorder is the stale "leading" entry, moving it
# to the end of the list (left-rotation by one position). Rows whose
# recorder list does not start with one of those leading entries are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

$recordedByCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Value2
    if ($header -eq "Recorded By") {
        $recordedByCol = $c
    }
}
if ($recordedByCol -eq 0) {
    $recordedByCol = 7
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val.Split(",")
        if ($parts.Length -gt 1) {
            $trimmedParts = @()
            foreach ($p in $parts) {
                $trimmedParts = $trimmedParts + $p.Trim()
            }

            $first = $trimmedParts[0]
            if ($first -eq "system" -or $first -eq "dnasr281@gmail.com") {
                $rest = $trimmedParts[1..($trimmedParts.Length - 1)]
                $rotatedParts = $rest + $first
                $newVal = $rotatedParts -join ", "
                $cell.Value2 = $newVal
            }
        }
    }
}
